# "finished cohort logic and updated sparse matrix"
#
# Renames the old "exclusions"/"adjusted" helper columns (G/H) to
# "exemptions"/"final", and appends a new "initial" column (N) with its
# cohort-inclusion flags, then refreshes the sheet/window chrome
# (selection, zoom, column widths) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2013 Cohort")

# --- header row -------------------------------------------------------
# Column order matters here: the shared-string table is rebuilt in
# first-use order, so touch G1, then the brand-new N1, then H1 so the
# new strings land as "exemptions", "initial", "final".
$ws.Range("G1").Value = "exemptions"
$ws.Range("N1").Value = "initial"
$ws.Range("H1").Value = "final"

# --- new "initial" column (N), rows 2-33 ------------------------------
$initialValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 0
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = 0
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = 0
    33 = 0
}

foreach ($r in $initialValues.Keys) {
    $ws.Cells.Item($r, 14).Value = $initialValues[$r]
}

# --- column widths -----------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 12.166666666666666
$ws.Columns.Item(13).ColumnWidth = 15.830729166666666

# --- selection / zoom ---------------------------------------------------
$ws.Range("I12").Select()
$excel.ActiveWindow.Zoom = 89
